$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "b857ce48-9fe8-42c3-83ee-14ba1287ffc3"
$ws.Range("B5").Value = "AMK"
$ws.Range("C5").Value = "Ang Mo Kio"
$ws.Range("D5").Value = "4"
